$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 306  # was 304
$ws.Range("F4").Value = 1247  # was 1238
$ws.Range("F5").Value = 357  # was 352
$ws.Range("F6").Value = 321  # was 317
$ws.Range("F7").Value = 3825  # was 3821
$ws.Range("F9").Value = 751  # was 744
$ws.Range("F10").Value = 1998  # was 1887
$ws.Range("F12").Value = 216  # was 215
$ws.Range("F13").Value = 730  # was 729
$ws.Range("F14").Value = 153  # was 149
$ws.Range("F15").Value = 154  # was 152
$ws.Range("F16").Value = 2084  # was 2081
$ws.Range("F18").Value = 6  # was 4
$ws.Range("F20").Value = 326  # was 324
$ws.Range("F21").Value = 223  # was 221
$ws.Range("F22").Value = 15  # was 14
$ws.Range("F23").Value = 267  # was 266

# Sheet 2: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 29  # was 28
$ws.Range("F7").Value = 32  # was 31
$ws.Range("F10").Value = 90  # was 89
$ws.Range("F12").Value = 82  # was 80
$ws.Range("F22").Value = 49  # was 48
$ws.Range("F23").Value = 56  # was 55

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 6392  # was 6390
$ws.Range("F3").Value = 811  # was 810
$ws.Range("F4").Value = 2071  # was 2068
$ws.Range("F5").Value = 304  # was 303

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 6392  # was 6390
$ws.Range("F3").Value = 811  # was 810
$ws.Range("F4").Value = 2071  # was 2068
$ws.Range("F5").Value = 304  # was 303
$ws.Range("F8").Value = 29  # was 28
$ws.Range("F12").Value = 306  # was 304
$ws.Range("F13").Value = 1247  # was 1238
$ws.Range("F14").Value = 357  # was 353
$ws.Range("F16").Value = 32  # was 31
$ws.Range("F18").Value = 321  # was 317
$ws.Range("F19").Value = 3825  # was 3821
$ws.Range("F22").Value = 90  # was 89
$ws.Range("F24").Value = 82  # was 80
$ws.Range("F25").Value = 751  # was 744
$ws.Range("F26").Value = 1998  # was 1887
$ws.Range("F29").Value = 216  # was 215
$ws.Range("F30").Value = 730  # was 729
$ws.Range("F31").Value = 153  # was 149
$ws.Range("F32").Value = 154  # was 152
$ws.Range("F34").Value = 2084  # was 2081
$ws.Range("F38").Value = 6  # was 4
$ws.Range("F40").Value = 326  # was 324
$ws.Range("F41").Value = 223  # was 221
$ws.Range("F42").Value = 15  # was 14
$ws.Range("F43").Value = 7  # was 6
$ws.Range("F48").Value = 49  # was 48
$ws.Range("F49").Value = 56  # was 55
$ws.Range("F50").Value = 267  # was 266

